$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Volume 30 Number 45" -> "Volume 30 Number 46" (A8) ---
$a8 = $ws.Range("A8")
$a8full = $a8.Value()
$idx = $a8full.LastIndexOf("45")
$a8.Characters($idx + 1, 2).Text = "46"

# --- Update week dates in C9: 11/6/2023 -> 11/13/2023, 11/12/2023 -> 11/19/2023 ---
$c9 = $ws.Range("C9")
$c9full = $c9.Value()
$idx1 = $c9full.IndexOf("11/6/2023")
$c9.Characters($idx1 + 1, 9).Text = "11/13/2023"
$c9full2 = $c9.Value()
$idx2 = $c9full2.IndexOf("11/12/2023")
$c9.Characters($idx2 + 1, 10).Text = "11/19/2023"

# --- Row 16-25, 27 numeric value updates (weekly crime statistics refresh) ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 97
$ws.Range("J16").Value = 101
$ws.Range("K16").Value = -3.960396039603
$ws.Range("L16").Value = 40.579710144927
$ws.Range("M16").Value = -51.010101010101
$ws.Range("N16").Value = -84.354838709677
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = -22.222222222222
$ws.Range("I17").Value = 142
$ws.Range("J17").Value = 150
$ws.Range("K17").Value = -5.333333333333
$ws.Range("L17").Value = -4.054054054054
$ws.Range("M17").Value = 24.561403508771
$ws.Range("N17").Value = -50.175438596491
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -83.333333333333
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -78.571428571428
$ws.Range("I18").Value = 80
$ws.Range("J18").Value = 107
$ws.Range("K18").Value = -25.233644859813
$ws.Range("L18").Value = -13.043478260869
$ws.Range("M18").Value = -68.253968253968
$ws.Range("N18").Value = -92.248062015503
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -13.333333333333
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 13.636363636363
$ws.Range("I19").Value = 581
$ws.Range("J19").Value = 465
$ws.Range("K19").Value = 24.946236559139
$ws.Range("L19").Value = 66.47564469914
$ws.Range("M19").Value = 29.111111111111
$ws.Range("N19").Value = 2.108963093145
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 14
$ws.Range("H20").Value = 16.666666666666
$ws.Range("I20").Value = 120
$ws.Range("J20").Value = 115
$ws.Range("K20").Value = 4.347826086956
$ws.Range("L20").Value = 62.162162162162
$ws.Range("M20").Value = -16.083916083916
$ws.Range("N20").Value = -95.100040832993
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -19.230769230769
$ws.Range("F21").Value = 81
$ws.Range("G21").Value = 85
$ws.Range("H21").Value = -4.705882352941
$ws.Range("I21").Value = 1035
$ws.Range("J21").Value = 947
$ws.Range("K21").Value = 9.292502639915
$ws.Range("L21").Value = 38.554216867469
$ws.Range("M21").Value = -11.914893617021
$ws.Range("N21").Value = -79.266826923076
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = 50
$ws.Range("I23").Value = 27
$ws.Range("K23").Value = -6.896551724137
$ws.Range("L23").Value = 12.5
$ws.Range("M23").Value = 8
$ws.Range("C24").Value = 29
$ws.Range("E24").Value = 26.086956521739
$ws.Range("G24").Value = 92
$ws.Range("H24").Value = 30.434782608695
$ws.Range("I24").Value = 1175
$ws.Range("J24").Value = 1016
$ws.Range("K24").Value = 15.649606299212
$ws.Range("L24").Value = 63.421418636995
$ws.Range("M24").Value = 29.834254143646
$ws.Range("C25").Value = 7
$ws.Range("E25").Value = 16.666666666666
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = 18.181818181818
$ws.Range("I25").Value = 242
$ws.Range("J25").Value = 211
$ws.Range("K25").Value = 14.691943127962
$ws.Range("L25").Value = 7.555555555555
$ws.Range("M25").Value = -23.899371069182
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50

# --- C23: change from text placeholder "0" to numeric value 2 (style matches F16) ---
$ws.Range("F16").Copy()
$ws.Range("C23").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C23").Value = 2

# --- C27: change from numeric value 2 to text placeholder "0" (style matches D27) ---
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("D27").Copy()
$ws.Range("C27").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
